$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header cells -------------------------------------------------
# A1 was "original test name" -> now "test location"
$ws.Range("A1").Value = "test location"

# New row 2: short note about the yellow-row convention used below
$ws.Range("A2").Value = "yellow rows: populated to RegressionTests"

# A3 was the bare "D:\...\Old\" path -> now annotated as historical/former location
$ws.Range("A3").Value = "Formerly: D:\PLanguage\PLang\Doc\Samples\New"

# A16 (second section title) was also the bare path -> now a different "Formerly:" note
$ws.Range("A16").Value = "Formerly: D:\PLanguage\PLang\Doc\Samples\New\Neg\Old\"

# --- Yellow highlight fill -------------------------------------------------
$yellow = 65535  # RGB(255,255,0)

# Row 6 (column header row for the first table): A,B,C,F,G get the yellow fill
foreach ($addr in @("A6", "B6", "C6", "F6", "G6")) {
    $ws.Range($addr).Interior.Color = $yellow
}

# --- Row 18 updates (regression test entry) --------------------------------
# A18: old test file path -> new repo-relative test location
$ws.Range("A18").Value = "Integration\DynamicError\Actions_1\Actions_1.p"

# D18/E18 used to be "?" placeholders; now populated
$ws.Range("D18").Value = "Yes"
$ws.Range("E18").Value = "No"

# Whole row 18 (A-F) gets the yellow "populated" highlight
foreach ($addr in @("A18", "B18", "C18", "D18", "E18", "F18")) {
    $ws.Range($addr).Interior.Color = $yellow
}

# --- Column width ------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 44.16666666666667

# --- Selection -----------------------------------------------------------
$ws.Rows("6:6").Select()
